$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row to append after existing data (row 83 -> row 84)
$newRow = 84

# Column A: date value, copy style/format from the cell above (A83)
$ws.Range("A83").Copy() | Out-Null
$ws.Range("A$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A$newRow").Value = 45884

# Column B: plain numeric value
$ws.Range("B$newRow").Value = 0

$excel.CutCopyMode = 0
